$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for the new "Core.Library.*" modules (ids 1000-1002) above the
#    existing "Cresus.Graph" (1003) row, and for more new modules (ids
#    1004-1012) between "Cresus.Graph" (1003) and the "Aider" (1013) row.
# ---------------------------------------------------------------------------
$ws.Rows("31:33").Insert()
$ws.Rows("38:40").Insert()

# After these inserts:
#   row 34 = old "Cresus.Graph" row (id 1003)            -> kept, restyled
#   row 35 = old "Cresus.Assets" row (id 1004)            -> repurposed
#   row 44 = old "Product.Aider" row (id 2000)            -> repurposed (id 1013)

# ---------------------------------------------------------------------------
# 2) Fill in the new rows. Values are entered in a specific order so that the
#    resulting shared-string table is built in the same sequence as the
#    reference workbook.
# ---------------------------------------------------------------------------

# id 1000
$ws.Range("A31").Value = 1000
$ws.Range("B31").Value = "Core.Library.Data"
$ws.Range("C31").Value = "s"

# id 1006
$ws.Range("A37").Value = 1006
$ws.Range("B37").Value = "Core.Library"
$ws.Range("C37").Value = "s"

# id 1007
$ws.Range("A38").Value = 1007
$ws.Range("B38").Value = "Core.Library.Address"
$ws.Range("C38").Value = "s"

# id 1009
$ws.Range("A40").Value = 1009
$ws.Range("B40").Value = "Core.Library.Documents"
$ws.Range("C40").Value = "s"

# id 1011
$ws.Range("A42").Value = 1011
$ws.Range("B42").Value = "Core.Library.Features"
$ws.Range("C42").Value = "s"

# id 1004 (replaces the old "Cresus.Assets" row)
$ws.Range("A35").Value = 1004
$ws.Range("B35").Value = "Core.Library.Finance"
$ws.Range("C35").Value = "s"
$ws.Range("D35").ClearContents()

# id 1001
$ws.Range("A32").Value = 1001
$ws.Range("B32").Value = "Core.Library.Images"
$ws.Range("C32").Value = "s"

# id 1010
$ws.Range("A41").Value = 1010
$ws.Range("B41").Value = "Core.Library.Measures"
$ws.Range("C41").Value = "s"

# id 1012
$ws.Range("A43").Value = 1012
$ws.Range("B43").Value = "Core.Library.UI"
$ws.Range("C43").Value = "s"

# id 1002
$ws.Range("A33").Value = 1002
$ws.Range("B33").Value = "Core.Library.UserManagement"
$ws.Range("C33").Value = "s"

# id 1005
$ws.Range("A36").Value = 1005
$ws.Range("B36").Value = "Core.Library.Workflows"
$ws.Range("C36").Value = "s"

# id 1008
$ws.Range("A39").Value = 1008
$ws.Range("B39").Value = "Core.Library.Business"
$ws.Range("C39").Value = "s"

# id 1014
$ws.Range("A45").Value = 1014
$ws.Range("B45").Value = "Data.Platform.SwissPostMatch"
$ws.Range("C45").Value = "A"

# id 1013 (replaces the old "Product.Aider" row, now at row 44)
$ws.Range("A44").Value = 1013
$ws.Range("B44").Value = "Aider"
$ws.Range("C44").Value = "A"
$ws.Range("D44").Value = "Produit ""AIDER"" développé pour l'EERV"

# id 2000 (new row, moved further down)
$ws.Range("A51").Value = 2000
$ws.Range("B51").Value = "Product.Assets.Data"
$ws.Range("C51").Value = "A"

# id 1016
$ws.Range("A47").Value = 1016
$ws.Range("B47").Value = "Cresus.Compta"
$ws.Range("C47").Value = "A"

# ---------------------------------------------------------------------------
# 3) Apply the italic style used to highlight the "Cresus.Graph" and "Aider"
#    rows.
# ---------------------------------------------------------------------------
$ws.Range("A34:D34").Font.Italic = $true
$ws.Range("A44:D44").Font.Italic = $true

# ---------------------------------------------------------------------------
# 4) Update the active selection to match the reference workbook.
# ---------------------------------------------------------------------------
$ws.Range("A47").Select()
